# Added Logout and closeBrowser code.
# Appends two new keyword rows ("logout", "closeBrowser") to the
# "FindAndBookFlight" sheet's keyword table (Table1423), growing it
# from A1:E6 to A1:E8, and leaves the selection on the last cell
# written (B8) to match the post-edit workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FindAndBookFlight")
$lo = $ws.ListObjects.Item("Table1423")

# Row 7: Sr No 6 -> logout
$lo.ListRows.Add() | Out-Null
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "logout"

# Row 8: Sr No 7 -> closeBrowser
$lo.ListRows.Add() | Out-Null
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "closeBrowser"

# Leave the selection on the last edited cell, as in the authored edit.
$ws.Range("B8").Select()
